$d = $word.ActiveDocument

# 1) Insert a new paragraph right after the "Nedan presenteras..." paragraph
#    with the "Vi förväntar oss..." text (near top of the document).
$pIntro = $d.Paragraphs.Item(3)
$pIntro.Range.InsertParagraphAfter()
$pNew = $d.Paragraphs.Item(4)
$pNew.Range.Text = "Vi förväntar oss att ni återkommer med ett skriftligt svar på vårt klagomål och även beskriver vilka korrigerande åtgärder ni satt in för att rätta till identifierade brister i er efterlevnad av den svenska FSC standarden."

# 2) Remove the two empty paragraphs and the duplicated "Vi förväntar oss..."
#    paragraph that used to sit at the very end of the document (after the
#    "...artskyddsförordningen" comment paragraph). Indices are +1 versus
#    the original document because of the paragraph inserted in step 1.
$pFirstToRemove = $d.Paragraphs.Item(38)
$pLastToRemove = $d.Paragraphs.Item(40)
$removeRange = $d.Range($pFirstToRemove.Range.Start, $pLastToRemove.Range.End)
$removeRange.Delete()

# 3) Update the date shown in the first-page header from 2023-11-13 to
#    2023-11-14.
$headers = $d.Sections.Item(1).Headers
$firstPageHeader = $headers.Item(2)
$firstPageHeader.Range.Find.Execute("2023-11-13", $true, $false, $false, $false, $false,
                                     $true, 1, $false, "2023-11-14", 2)
